# Insert a new weekly price record as row 80 in the "Brócoli" sheet,
# pushing the existing rows 80:200 down to 81:201 (dimension grows to R201).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 80 downward by inserting a new blank row at position 80.
$ws.Rows.Item(80).Insert()

# Populate the newly inserted row 80 with the new record's data.
$ws.Cells.Item(80, 1).Value = 4
$ws.Cells.Item(80, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(80, 3).Value = 'Los Lagos'
$ws.Cells.Item(80, 4).Value = 44482
$ws.Cells.Item(80, 5).Value = 10
$ws.Cells.Item(80, 6).Value = 100112023
$ws.Cells.Item(80, 7).Value = 'Brócoli'
$ws.Cells.Item(80, 8).Value = 'Sin especificar'
$ws.Cells.Item(80, 9).Value = 'Segunda'
$ws.Cells.Item(80, 10).Value = 250
$ws.Cells.Item(80, 11).Value = 1000
$ws.Cells.Item(80, 12).Value = 1000
$ws.Cells.Item(80, 13).Value = 1000
$ws.Cells.Item(80, 14).Value = '$/unidad'
$ws.Cells.Item(80, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(80, 16).Value = 1000
$ws.Cells.Item(80, 17).Value = 1
$ws.Cells.Item(80, 18).Value = 'Hortaliza'
